# The sheet contains a weekly time-series of "Zapallo italiano" price records for
# "Macroferia Regional de Talca" (rows 248..313). A new weekly record was inserted
# at the top of that series (new row 248), which pushes every existing record down
# by one row; the former last record (old row 313) becomes the new last record
# (row 314). Columns A,B,C,E,F,G,H,I,R (market/category metadata) are identical
# for every record in this block, while D (Fecha), J (Volumen), K/L/M (precios),
# N (Unidad), O (Origen), P (Precio $/Kg) and Q (Kg o Unidades) carry the
# per-record data that shifts down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Create the new last row (314) by duplicating the constant metadata columns
#    from the (current) last row of the series, row 313.
$srcRow = 313
$dstRow = 314
$ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2   # Mercado ID
$ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2   # Mercado
$ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2   # Región
$ws.Cells.Item($dstRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value2   # Codreg
$ws.Cells.Item($dstRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value2   # Categoría ID
$ws.Cells.Item($dstRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value2   # Categoría
$ws.Cells.Item($dstRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value2   # Variedad
$ws.Cells.Item($dstRow, 9).Value = $ws.Cells.Item($srcRow, 9).Value2   # Calidad
$ws.Cells.Item($dstRow, 18).Value = $ws.Cells.Item($srcRow, 18).Value2 # Clasificación
$ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat

# 2) Row 314 gets the data values that used to belong to row 313 (shifted down).
$ws.Cells.Item(314, 4).Value = 44544
$ws.Cells.Item(314, 10).Value = 400
$ws.Cells.Item(314, 11).Value = 4500
$ws.Cells.Item(314, 12).Value = 4500
$ws.Cells.Item(314, 13).Value = 4500
$ws.Cells.Item(314, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(314, 15).Value = "Región del Maule"
$ws.Cells.Item(314, 16).Value = 75
$ws.Cells.Item(314, 17).Value = 60

# 3) Shift the data values of rows 249..313 down by one: new row r takes on
#    what used to be in row (r-1). Walk from the bottom up so a row's old
#    value is always read before it gets overwritten.
for ($r = 313; $r -ge 249; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($src, 14).Value2
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($src, 15).Value2
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($src, 16).Value2
    $ws.Cells.Item($r, 17).Value = $ws.Cells.Item($src, 17).Value2
}

# 4) Row 248 finally receives the brand-new record's data.
$ws.Cells.Item(248, 4).Value = 44642
$ws.Cells.Item(248, 10).Value = 400
$ws.Cells.Item(248, 11).Value = 6000
$ws.Cells.Item(248, 12).Value = 6000
$ws.Cells.Item(248, 13).Value = 6000
$ws.Cells.Item(248, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(248, 15).Value = "Región del Maule"
$ws.Cells.Item(248, 16).Value = 120
$ws.Cells.Item(248, 17).Value = 50

$addr = $ws.UsedRange.Address()
Write-Host "UsedRange: $addr"
